$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 71.947365
$ws.Range("I2").Value = 19.875
$ws.Range("K2").Value = 19.875
$ws.Range("M2").Value = 93.125
$ws.Range("H28").Value = 800
$ws.Range("J28").Value = 1000
$ws.Range("L28").Value = 1000
$ws.Range("N28").Value = -1970
$ws.Range("H33").Value = 1771.091
$ws.Range("I33").Value = 2061.5
$ws.Range("K33").Value = 2061.5
$ws.Range("M33").Value = -1832.5
$ws.Range("H43").Value = 456458.22
$ws.Range("J43").Value = 456458.22
$ws.Range("L43").Value = 456458.22
$ws.Range("N43").Value = -456596.22
$ws.Range("H51").Value = 12229.462
$ws.Range("J51").Value = 8225.817999999999
$ws.Range("L51").Value = 8225.817999999999
$ws.Range("N51").Value = -9193.817999999999
$ws.Range("H64").Value = 4631.8335
$ws.Range("I64").Value = 4325.25
$ws.Range("J64").Value = 5245
$ws.Range("K64").Value = 4325.25
$ws.Range("L64").Value = 5245
$ws.Range("M64").Value = -4077.25
$ws.Range("N64").Value = -5741
$ws.Range("H67").Value = 4631.8335
$ws.Range("I67").Value = 4325.25
$ws.Range("J67").Value = 5245
$ws.Range("K67").Value = 4325.25
$ws.Range("L67").Value = 5245
$ws.Range("M67").Value = -3467.25
$ws.Range("N67").Value = -6961
$ws.Range("H69").Value = 100000
$ws.Range("J69").Value = 100000
$ws.Range("L69").Value = 300000
$ws.Range("N69").Value = -301748
$ws.Range("H72").Value = 100000
$ws.Range("J72").Value = 100000
$ws.Range("L72").Value = 900000
$ws.Range("N72").Value = -908736
$ws.Range("H127").Value = 3346.2727
$ws.Range("I127").Value = 3346.2727
$ws.Range("K127").Value = 10038.8181
$ws.Range("M127").Value = -5078.8181
$ws.Range("H129").Value = 1292.8572
$ws.Range("I129").Value = 762.5
$ws.Range("K129").Value = 2287.5
$ws.Range("M129").Value = 2712.5
$ws.Range("H135").Value = 1251268.5
$ws.Range("I135").Value = 1429750
$ws.Range("K135").Value = 12867750
$ws.Range("M135").Value = -12865215
$ws.Range("H138").Value = 2446.0322
$ws.Range("J138").Value = 2441.25
$ws.Range("L138").Value = 7323.75
$ws.Range("N138").Value = -17603.75
$ws.Range("H141").Value = 2854.7778
$ws.Range("I141").Value = 2854.7778
$ws.Range("K141").Value = 8564.3334
$ws.Range("M141").Value = -3384.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1626979.9
$ws.Range("I32").Value = 1626979.9
$ws.Range("K32").Value = 1626979.9
$ws.Range("M32").Value = -1626692.9
$ws.Range("H61").Value = 4052.9355
$ws.Range("I61").Value = 1773.1459
$ws.Range("J61").Value = 11869.357
$ws.Range("K61").Value = 1773.1459
$ws.Range("L61").Value = 11869.357
$ws.Range("M61").Value = -1561.1459
$ws.Range("N61").Value = -12293.357
$ws.Range("H122").Value = 16124.1875
$ws.Range("I122").Value = 21498.7
$ws.Range("J122").Value = 7166.6665
$ws.Range("K122").Value = 64496.10000000001
$ws.Range("L122").Value = 21499.9995
$ws.Range("M122").Value = -62046.10000000001
$ws.Range("N122").Value = -26399.9995
$ws.Range("H136").Value = 4052.9355
$ws.Range("I136").Value = 1773.1459
$ws.Range("J136").Value = 11869.357
$ws.Range("K136").Value = 5319.4377
$ws.Range("L136").Value = 35608.071
$ws.Range("M136").Value = -2769.4377
$ws.Range("N136").Value = -40708.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 2569.2856
$ws.Range("I128").Value = 2569.2856
$ws.Range("K128").Value = 7707.8568
$ws.Range("M128").Value = -5217.8568
$ws.Range("H134").Value = 6003.61
$ws.Range("I134").Value = 1881.7059
$ws.Range("K134").Value = 5645.1177
$ws.Range("M134").Value = -3110.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3012.8235
$ws.Range("I16").Value = 1893.1428
$ws.Range("J16").Value = 4821.5386
$ws.Range("K16").Value = 1893.1428
$ws.Range("L16").Value = 4821.5386
$ws.Range("M16").Value = -1606.1428
$ws.Range("N16").Value = -5395.5386
$ws.Range("H22").Value = 297.375
$ws.Range("I22").Value = 363.75
$ws.Range("J22").Value = 231
$ws.Range("K22").Value = 363.75
$ws.Range("L22").Value = 231
$ws.Range("M22").Value = -13.75
$ws.Range("N22").Value = -931
$ws.Range("H31").Value = 7414001.5
$ws.Range("I31").Value = 2241.0386
$ws.Range("J31").Value = 17556410
$ws.Range("K31").Value = 2241.0386
$ws.Range("L31").Value = 17556410
$ws.Range("M31").Value = -1946.0386
$ws.Range("N31").Value = -17557000
$ws.Range("H34").Value = 7414001.5
$ws.Range("I34").Value = 2241.0386
$ws.Range("J34").Value = 17556410
$ws.Range("K34").Value = 2241.0386
$ws.Range("L34").Value = 17556410
$ws.Range("M34").Value = -2039.0386
$ws.Range("N34").Value = -17556814
$ws.Range("H58").Value = 7843.069
$ws.Range("I58").Value = 2109.5
$ws.Range("K58").Value = 2109.5
$ws.Range("M58").Value = -1906.5
$ws.Range("H99").Value = 11809.667
$ws.Range("I99").Value = 13258
$ws.Range("K99").Value = 13258
$ws.Range("M99").Value = -11760
$ws.Range("H107").Value = 2453.375
$ws.Range("J107").Value = 2848
$ws.Range("L107").Value = 2848
$ws.Range("N107").Value = -6688
$ws.Range("H113").Value = 3012.8235
$ws.Range("I113").Value = 1893.1428
$ws.Range("J113").Value = 4821.5386
$ws.Range("K113").Value = 1893.1428
$ws.Range("L113").Value = 4821.5386
$ws.Range("M113").Value = 276.8571999999999
$ws.Range("N113").Value = -9161.5386
$ws.Range("H126").Value = 11809.667
$ws.Range("I126").Value = 13258
$ws.Range("K126").Value = 39774
$ws.Range("M126").Value = -37304
$ws.Range("H132").Value = 10816342
$ws.Range("I132").Value = 2360.35
$ws.Range("K132").Value = 7081.049999999999
$ws.Range("M132").Value = -4551.049999999999
$ws.Range("H136").Value = 7843.069
$ws.Range("I136").Value = 2109.5
$ws.Range("K136").Value = 6328.5
$ws.Range("M136").Value = -3778.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83801.25
$ws.Range("I2").Value = 64.625
$ws.Range("K2").Value = 387.75
$ws.Range("M2").Value = -274.75
$ws.Range("H92").Value = 1440.7333
$ws.Range("I92").Value = 1339.6666
$ws.Range("J92").Value = 1466
$ws.Range("K92").Value = 4018.9998
$ws.Range("L92").Value = 4398
$ws.Range("M92").Value = -2770.9998
$ws.Range("N92").Value = -6894
$ws.Range("H109").Value = 992.6667
$ws.Range("I109").Value = 992.6667
$ws.Range("K109").Value = 2978.0001
$ws.Range("M109").Value = -1938.0001
$ws.Range("H122").Value = 5659405.5
$ws.Range("J122").Value = 3750
$ws.Range("L122").Value = 33750
$ws.Range("N122").Value = -38650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6022
$ws.Range("I113").Value = 3075.15
$ws.Range("J113").Value = 9296.277
$ws.Range("K113").Value = 3075.15
$ws.Range("L113").Value = 9296.277
$ws.Range("M113").Value = -905.1500000000001
$ws.Range("N113").Value = -13636.277
$ws.Range("H132").Value = 4355.9585
$ws.Range("I132").Value = 1794.5
$ws.Range("J132").Value = 9478.875
$ws.Range("K132").Value = 5383.5
$ws.Range("L132").Value = 28436.625
$ws.Range("M132").Value = -2853.5
$ws.Range("N132").Value = -33496.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1186.8064
$ws.Range("I22").Value = 679.6
$ws.Range("J22").Value = 3300.1667
$ws.Range("K22").Value = 679.6
$ws.Range("L22").Value = 3300.1667
$ws.Range("M22").Value = -384.6
$ws.Range("N22").Value = -3890.1667
$ws.Range("H27").Value = 1186.8064
$ws.Range("I27").Value = 679.6
$ws.Range("J27").Value = 3300.1667
$ws.Range("K27").Value = 679.6
$ws.Range("L27").Value = 3300.1667
$ws.Range("M27").Value = -572.6
$ws.Range("N27").Value = -3514.1667
$ws.Range("H40").Value = 7217.483
$ws.Range("I40").Value = 6328.857
$ws.Range("K40").Value = 6328.857
$ws.Range("M40").Value = -6192.857
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H46").Value = 1380790.8
$ws.Range("J46").Value = 2285.7144
$ws.Range("L46").Value = 2285.7144
$ws.Range("N46").Value = -2661.7144
$ws.Range("H55").Value = 335.33334
$ws.Range("I55").Value = 89
$ws.Range("J55").Value = 581.6667
$ws.Range("K55").Value = 89
$ws.Range("L55").Value = 581.6667
$ws.Range("M55").Value = 84
$ws.Range("N55").Value = -927.6667
$ws.Range("H132").Value = 6594.375
$ws.Range("J132").Value = 9407.73
$ws.Range("L132").Value = 28223.19
$ws.Range("N132").Value = -33283.19

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1942.1111
$ws.Range("I113").Value = 1257.8334
$ws.Range("K113").Value = 3773.5002
$ws.Range("M113").Value = -1603.5002
$ws.Range("H122").Value = 141876.34
$ws.Range("I122").Value = 237416.4
$ws.Range("J122").Value = 6527.9165
$ws.Range("K122").Value = 712249.2
$ws.Range("L122").Value = 19583.7495
$ws.Range("M122").Value = -709799.2
$ws.Range("N122").Value = -24483.7495
$ws.Range("H132").Value = 10750.923
$ws.Range("I132").Value = 50252
$ws.Range("J132").Value = 3568.9092
$ws.Range("K132").Value = 150756
$ws.Range("L132").Value = 10706.7276
$ws.Range("M132").Value = -148226
$ws.Range("N132").Value = -15766.7276
$ws.Range("H136").Value = 43049.926
$ws.Range("I136").Value = 1596.6666
$ws.Range("K136").Value = 4789.9998
$ws.Range("M136").Value = -2239.9998
$ws.Range("H139").Value = 92662
$ws.Range("J139").Value = 93549.336
$ws.Range("L139").Value = 93549.336
